$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the old grouped-header cells in row 1 ---
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# --- Rewrite row 1 as a flattened header row (was a 2-row MultiIndex header) ---
$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# --- Row 2 keeps the original detailed sub-header text, but is now hidden ---
$ws.Rows(2).Hidden = $true

# --- New blank separator row 3, also hidden ---
$ws.Rows(3).Hidden = $true

# --- Fill in the previously-missing "Tkl%" = 0 cells for a handful of rows ---
$ws.Range("O7").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("O19").Value = 0

# --- Totals row is now hidden ---
$ws.Rows(20).Hidden = $true

# --- Match the author's final cursor position ---
$ws.Range("O21").Select()
